$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Collapse the proofErr-split runs ("Using <spellcheck>DxDesigner</spellcheck>
#    with PADS" etc.) into single, plain runs with no spell-check wrapper.
#    A same-text Find & Replace rewrites the matched range as one run, which
#    drops the w:proofErr markers and merges the surrounding runs.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("Using DxDesigner with PADS", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Using DxDesigner with PADS", 2) | Out-Null

$d.Content.Find.Execute("Creating DxDesigner parts", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Creating DxDesigner parts", 2) | Out-Null

$d.Content.Find.Execute("Editing DxDesigner parts", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Editing DxDesigner parts", 2) | Out-Null

$d.Content.Find.Execute("Blind and buried vias", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Blind and buried vias", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Design choices?" -> "Design choices" (drop the trailing question mark).
#    A placeholder "X" is appended so we have real content to anchor a Find
#    on afterwards (collapsed zero-length ranges built straight from
#    Range(start,start) land at the wrong spot in this host, but a Find
#    result collapsed to its Start lands correctly).
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("Design choices?", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Design choicesX", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Relocate the _GoBack bookmark from the end of "EMI considerations" to
#    the end of "Design choices".
# ---------------------------------------------------------------------------

$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

$xr = $d.Content
$xr.Find.Execute("X", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$xr.Collapse(1)
$d.Bookmarks.Add("_GoBack", $xr)

# Remove the placeholder "X" now that the bookmark is anchored in place.
$xr2 = $d.Content
$xr2.Find.Execute("X", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$xr2.Text = ""
